$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "26.398.87"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "1.622.56"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.42"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.496"
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0621"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.91"
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0840"
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").Value = "1.848.80"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "1.635.65"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.10"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.04"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D17").Value = "26.392.83"
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").Value = "0.0₃0738"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.62"
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.21"
$ws.Range("E22").Value = "  +2.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.28"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("E24").Value = "  +3.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.51"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.82"
$ws.Range("E28").Value = "  +2.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.55"
$ws.Range("E29").Value = "  +1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0507"
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("D36").Value = "1.210.97"
$ws.Range("E36").Value = "  +3.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0172"
$ws.Range("E37").Value = "  +3.19%  "
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.797"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.500"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("E41").Value = "  -3.25%  "
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").Value = "1.759.69"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.60"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.58"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.58"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.61"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.407"
$ws.Range("E51").Value = "  -0.52%  "
